# Generate Report for Handoff
# Updates the localization-status report: status moves from "In Translation"
# to "Ready for handoff" and the associated timestamps are refreshed, on all
# three sheets (Overview, zh-cn, de-de). Also widens the "Status" columns
# (E/F on Overview, C on the per-locale sheets) to fit the new text.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet --------------------------------------------------
$ws = $wb.Sheets.Item("Overview")
$ws.Range("E2").Value = "Ready for handoff"
$ws.Range("F2").Value = "Ready for handoff"
$ws.Range("G2").Value = "2016-08-17 22:38:07"
$ws.Columns.Item(5).ColumnWidth = 16.38
$ws.Columns.Item(6).ColumnWidth = 16.38

# ---- zh-cn sheet -------------------------------------------------------
$ws = $wb.Sheets.Item("zh-cn")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("H2").Value = "2016-08-17 22:37:57"
$ws.Columns.Item(3).ColumnWidth = 16.38

# ---- de-de sheet ---------------------------------------------------
$ws = $wb.Sheets.Item("de-de")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("H2").Value = "2016-08-17 22:38:07"
$ws.Columns.Item(3).ColumnWidth = 16.38
